$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.033.40"
$ws.Range("E2").Value = "  +1.80%  "

$ws.Range("D3").Value = "2.587.97"
$ws.Range("E3").Value = "  +0.51%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "522.90"
$ws.Range("E5").Value = "  +0.95%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.22"
$ws.Range("E6").Value = "  -1.98%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("D9").Value = "2.599.99"
$ws.Range("E9").Value = "  +0.43%  "

$ws.Range("E10").Value = "  -2.95%  "

$ws.Range("E11").Value = "  -0.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.329"
$ws.Range("E12").Value = "  +1.20%  "

$ws.Range("E13").Value = "  +3.32%  "

$ws.Range("D14").Value = "3.043.40"
$ws.Range("E14").Value = "  +0.56%  "

$ws.Range("D15").Value = "58.975.91"
$ws.Range("E15").Value = "  +1.75%  "

$ws.Range("E16").Value = "  +0.60%  "

$ws.Range("D17").Value = "2.586.87"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("E18").Value = "  -0.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.54"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("E20").Value = "  -0.03%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.08"
$ws.Range("E21").Value = "  -1.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.50"
$ws.Range("E22").Value = "  +2.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.08"
$ws.Range("E24").Value = "  +0.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.168"
$ws.Range("E25").Value = "  +1.51%  "

$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("E27").Value = "  +0.07%  "

$ws.Range("E28").Value = "  +0.48%  "

$ws.Range("E29").Value = "  +0.00%  "

$ws.Range("D30").Value = "0.0₃0725"
$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.91"
$ws.Range("E31").Value = "  -4.88%  "

$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.69"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("E34").Value = "  -0.58%  "

$ws.Range("E35").Value = "  -0.51%  "

$ws.Range("E36").Value = "  -1.68%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "36.81"
$ws.Range("E37").Value = "  +2.35%  "

$ws.Range("E38").Value = "  +1.53%  "

$ws.Range("E39").Value = "  -0.91%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.816"
$ws.Range("E40").Value = "  -5.74%  "

$ws.Range("E41").Value = "  -0.45%  "

$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "272.12"
$ws.Range("E43").Value = "  +0.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.77"
$ws.Range("E44").Value = "  +0.95%  "

$ws.Range("E45").Value = "  +0.36%  "

$ws.Range("E46").Value = "  +0.47%  "

$ws.Range("E47").Value = "  -1.06%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.42"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("D49").Value = "1.965.33"
$ws.Range("E49").Value = "  -0.25%  "

$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("E51").Value = "  -0.18%  "
